$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to hold literal TEXT even when the string looks like
# a number or a date (Excel would otherwise auto-convert "2024-03-01" into a
# date serial, or "1" into the number 1). Switching NumberFormat to Text
# before the assignment keeps the literal string; switching the Style back
# to "Normal" afterwards removes the temporary formatting again so the cell
# keeps the workbook's default (unstyled) look.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 7 ("jalebi"): quantity/price were stored as text ("1"/"20"); convert
# them to real numbers. expiry (E7) already reads "2024-03-12" - no change.
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 20

# Row 8 (new): siddhu / 2024-03-01 / 1 / 2 / 2024-03-14
$ws.Range("A8").Value = "siddhu"
Set-TextValue $ws.Range("B8") "2024-03-01"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 2
Set-TextValue $ws.Range("E8") "2024-03-14"

# Row 9 (new): siddharth / 2024-03-14 / "1" / "2" / 2024-03-14
# Here quantity and price are stored as TEXT ("1"/"2"), not numbers.
$ws.Range("A9").Value = "siddharth"
Set-TextValue $ws.Range("B9") "2024-03-14"
Set-TextValue $ws.Range("C9") "1"
Set-TextValue $ws.Range("D9") "2"
Set-TextValue $ws.Range("E9") "2024-03-14"
